$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells in column D whose new values look numeric need to be pinned
# to Text format first so Excel keeps them as literal strings (matching the
# source data, e.g. "74.46") instead of silently parsing them into numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.451.68"
$ws.Range("D3").Value = "2.234.10"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "244.34"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "74.46"
$ws.Range("E7").Value = "  -4.18%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").Value = "43.21"
$ws.Range("E10").Value = "  +2.14%  "
$ws.Range("D11").Value = "0.0945"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").Value = "7.11"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "14.42"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("D15").Value = "0.850"
$ws.Range("E15").Value = "  -1.49%  "
$ws.Range("D16").Value = "2.242.34"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "42.197.73"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "0.0000103"
$ws.Range("E18").Value = "  +4.99%  "
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "71.86"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "10.19"
$ws.Range("E21").Value = "  +39.27%  "
$ws.Range("D22").Value = "230.33"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("E23").Value = "  -6.73%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "11.56"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("D26").Value = "3.66"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +3.93%  "
$ws.Range("D29").Value = "166.37"
$ws.Range("E29").Value = "  -2.05%  "
$ws.Range("D30").Value = "20.87"
$ws.Range("E30").Value = "  +1.33%  "
$ws.Range("D31").Value = "5.86"
$ws.Range("E31").Value = "  +19.48%  "
$ws.Range("D32").Value = "0.0805"
$ws.Range("E32").Value = "  -3.29%  "
$ws.Range("E33").Value = "  -3.36%  "
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "29.51"
$ws.Range("E35").Value = "  -12.59%  "
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("D38").Value = "13.20"
$ws.Range("E38").Value = "  -7.93%  "
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("D40").Value = "5.67"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").Value = "63.11"
$ws.Range("D42").Value = "0.199"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "8.80"
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("D44").Value = "105.61"
$ws.Range("E44").Value = "  -5.94%  "
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  +3.79%  "
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  -0.38%  "
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "4.09"
$ws.Range("E51").Value = "  -2.77%  "
